$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue([string]$cellRef, [string]$value) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue "D2" "41.948.73"
$ws.Range("E2").Value = "  -0.60%  "
Set-TextValue "D3" "2.214.38"
$ws.Range("E3").Value = "  -1.30%  "
$ws.Range("E4").Value = "  +0.01%  "
Set-TextValue "D5" "240.92"
$ws.Range("E5").Value = "  -2.21%  "
Set-TextValue "D6" "0.616"
$ws.Range("E6").Value = "  -2.07%  "
Set-TextValue "D7" "73.11"
$ws.Range("E7").Value = "  -1.71%  "
$ws.Range("E8").Value = "  -0.03%  "
Set-TextValue "D9" "0.605"
$ws.Range("E9").Value = "  -1.96%  "
Set-TextValue "D10" "42.87"
$ws.Range("E10").Value = "  +1.64%  "
Set-TextValue "D11" "0.0952"
$ws.Range("E11").Value = "  +1.00%  "
Set-TextValue "D12" "7.08"
$ws.Range("E12").Value = "  -1.23%  "
$ws.Range("E13").Value = "  -0.37%  "
Set-TextValue "D14" "2.545.72"
$ws.Range("E14").Value = "  -1.31%  "
Set-TextValue "D15" "14.21"
$ws.Range("E15").Value = "  -2.06%  "
Set-TextValue "D16" "0.835"
$ws.Range("E16").Value = "  -2.04%  "
Set-TextValue "D17" "2.207.26"
$ws.Range("E17").Value = "  -1.05%  "
Set-TextValue "D18" "41.784.12"
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("E19").Value = "  +10.53%  "
Set-TextValue "D20" "73.15"
$ws.Range("E20").Value = "  +1.37%  "
Set-TextValue "D21" "6.15"
$ws.Range("E21").Value = "  +0.26%  "
Set-TextValue "D22" "10.29"
$ws.Range("E22").Value = "  +17.46%  "
Set-TextValue "D23" "228.77"
$ws.Range("E23").Value = "  -1.12%  "
Set-TextValue "D24" "2.09"
$ws.Range("E24").Value = "  -5.61%  "
Set-TextValue "D25" "11.64"
$ws.Range("E25").Value = "  +1.00%  "
Set-TextValue "D26" "1.00"
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("E28").Value = "  -1.83%  "
Set-TextValue "D29" "2.19"
$ws.Range("E29").Value = "  +1.13%  "
Set-TextValue "D30" "166.57"
$ws.Range("E30").Value = "  -1.45%  "
Set-TextValue "D31" "20.58"
$ws.Range("E31").Value = "  -0.35%  "
Set-TextValue "D32" "5.68"
$ws.Range("E32").Value = "  +9.17%  "
Set-TextValue "D34" "0.125"
$ws.Range("E34").Value = "  -0.30%  "
Set-TextValue "D35" "29.12"
$ws.Range("E35").Value = "  -7.11%  "
$ws.Range("E36").Value = "  -7.98%  "
$ws.Range("E37").Value = "  -5.19%  "
Set-TextValue "D38" "0.0300"
$ws.Range("E38").Value = "  -4.28%  "
Set-TextValue "D39" "13.48"
$ws.Range("E39").Value = "  -2.12%  "
Set-TextValue "D40" "66.23"
$ws.Range("E40").Value = "  +6.52%  "
$ws.Range("E41").Value = "  -3.19%  "
$ws.Range("E42").Value = "  -2.58%  "
$ws.Range("E43").Value = "  -3.95%  "
Set-TextValue "D44" "8.66"
$ws.Range("E44").Value = "  +0.03%  "
Set-TextValue "D45" "103.55"
$ws.Range("E45").Value = "  -2.83%  "
$ws.Range("E46").Value = "  -2.34%  "
Set-TextValue "D47" "2.38"
$ws.Range("E47").Value = "  +4.46%  "
Set-TextValue "D48" "1.12"
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("E50").Value = "  -0.37%  "
Set-TextValue "D51" "2.420.55"
$ws.Range("E51").Value = "  -1.33%  "

Write-Host "Applied crypto list updates"
